# Updated cryptos list with GitHub Actions.
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) snapshot for
# each coin row, and re-sorts rows 15/16 (Chainlink now outranks WrappedEther).
#
# Column D values are written with a leading apostrophe so Excel keeps
# numeric-looking prices (e.g. 98.81) as literal text instead of silently
# coercing them into floating-point numbers (matching the sheet's existing
# text-formatted Price column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'42.528.24"
$ws.Cells.Item(2, 5).Value = "  -0.62%  "
$ws.Cells.Item(3, 4).Value = "'2.521.31"
$ws.Cells.Item(3, 5).Value = "  -1.09%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'312.12"
$ws.Cells.Item(5, 5).Value = "  +1.12%  "
$ws.Cells.Item(6, 4).Value = "'98.81"
$ws.Cells.Item(6, 5).Value = "  -2.52%  "
$ws.Cells.Item(7, 5).Value = "  -1.28%  "
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
$ws.Cells.Item(9, 4).Value = "'0.518"
$ws.Cells.Item(9, 5).Value = "  -2.89%  "
$ws.Cells.Item(10, 4).Value = "'35.35"
$ws.Cells.Item(10, 5).Value = "  -2.44%  "
$ws.Cells.Item(11, 4).Value = "'0.0801"
$ws.Cells.Item(11, 5).Value = "  -0.92%  "
$ws.Cells.Item(12, 5).Value = "  +0.53%  "
$ws.Cells.Item(13, 4).Value = "'7.23"
$ws.Cells.Item(13, 5).Value = "  -2.44%  "
$ws.Cells.Item(14, 4).Value = "'2.910.75"
$ws.Cells.Item(14, 5).Value = "  -0.96%  "
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "'15.28"
$ws.Cells.Item(15, 5).Value = "  -4.35%  "
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "'2.516.43"
$ws.Cells.Item(16, 5).Value = "  -3.29%  "
$ws.Cells.Item(17, 4).Value = "'0.808"
$ws.Cells.Item(17, 5).Value = "  -3.83%  "
$ws.Cells.Item(18, 4).Value = "'42.542.66"
$ws.Cells.Item(18, 5).Value = "  -0.67%  "
$ws.Cells.Item(19, 4).Value = "'6.59"
$ws.Cells.Item(19, 5).Value = "  -2.75%  "
$ws.Cells.Item(20, 5).Value = "  -0.93%  "
$ws.Cells.Item(21, 4).Value = "'12.17"
$ws.Cells.Item(21, 5).Value = "  -1.57%  "
$ws.Cells.Item(22, 4).Value = "'69.38"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "
$ws.Cells.Item(23, 4).Value = "'241.35"
$ws.Cells.Item(23, 5).Value = "  -2.86%  "
$ws.Cells.Item(24, 5).Value = "  -1.38%  "
$ws.Cells.Item(25, 5).Value = "  -3.39%  "
$ws.Cells.Item(27, 4).Value = "'25.45"
$ws.Cells.Item(27, 5).Value = "  -4.15%  "
$ws.Cells.Item(28, 4).Value = "'2.26"
$ws.Cells.Item(28, 5).Value = "  -4.02%  "
$ws.Cells.Item(29, 4).Value = "'10.07"
$ws.Cells.Item(29, 5).Value = "  -0.65%  "
$ws.Cells.Item(30, 4).Value = "'38.22"
$ws.Cells.Item(30, 5).Value = "  -5.86%  "
$ws.Cells.Item(31, 4).Value = "'5.82"
$ws.Cells.Item(31, 5).Value = "  +1.49%  "
$ws.Cells.Item(32, 4).Value = "'156.98"
$ws.Cells.Item(32, 5).Value = "  +0.24%  "
$ws.Cells.Item(33, 5).Value = "  +1.23%  "
$ws.Cells.Item(34, 4).Value = "'2.66"
$ws.Cells.Item(34, 5).Value = "  +1.84%  "
$ws.Cells.Item(35, 5).Value = "  -2.20%  "
$ws.Cells.Item(36, 4).Value = "'3.15"
$ws.Cells.Item(36, 5).Value = "  -4.67%  "
$ws.Cells.Item(37, 5).Value = "  -6.67%  "
$ws.Cells.Item(38, 4).Value = "'17.62"
$ws.Cells.Item(38, 5).Value = "  -4.02%  "
$ws.Cells.Item(39, 5).Value = "  -1.68%  "
$ws.Cells.Item(40, 5).Value = "  -0.85%  "
$ws.Cells.Item(41, 4).Value = "'4.13"
$ws.Cells.Item(41, 5).Value = "  -2.56%  "
$ws.Cells.Item(42, 4).Value = "'21.82"
$ws.Cells.Item(42, 5).Value = "  -3.45%  "
$ws.Cells.Item(43, 5).Value = "  +0.20%  "
$ws.Cells.Item(44, 5).Value = "  -0.24%  "
$ws.Cells.Item(45, 5).Value = "  -1.19%  "
$ws.Cells.Item(46, 4).Value = "'1.988.76"
$ws.Cells.Item(46, 5).Value = "  +0.10%  "
$ws.Cells.Item(47, 4).Value = "'9.02"
$ws.Cells.Item(47, 5).Value = "  -0.02%  "
$ws.Cells.Item(48, 4).Value = "'2.772.44"
$ws.Cells.Item(48, 5).Value = "  -0.69%  "
$ws.Cells.Item(49, 5).Value = "  -2.28%  "
$ws.Cells.Item(50, 4).Value = "'78.76"
$ws.Cells.Item(50, 5).Value = "  -3.23%  "
$ws.Cells.Item(51, 4).Value = "'71.56"
$ws.Cells.Item(51, 5).Value = "  -2.82%  "
